# Generate Report for Handoff
# Renamed source file 759239cf-2818-4f34-9c84-0fae4df38b1c.md -> f425f407-2121-4432-b119-c1a0d4ae2d42.md
# and refreshed the localization-status report (new handoff xliff files / timestamps,
# cleared stale handback info pending the new cycle).

$wb = $excel.ActiveWorkbook

$oldId = "759239cf-2818-4f34-9c84-0fae4df38b1c"
$newId = "f425f407-2121-4432-b119-c1a0d4ae2d42"
$newHash = "5e20e21316651e5fa0b3c13e7a6aec1a1732c791"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newId.md"
$wsOverview.Range("G2").Value = "2016-08-19 15:06:36"

$overviewLinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2978d00b5a530dcd104d2590d7e7202a99367c54/e2e/$oldId.md"
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $overviewLinkUrl, "", "", "e2e\$newId.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("G2").Value = "$newId.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-19 15:06:32"
$wsZhCn.Range("I2").Value = ""
$wsZhCn.Range("I2").Style = "Normal"
$wsZhCn.Range("J2").Value = ""
$wsZhCn.Range("K2").Value = "0001-01-01 00:00:00"

$zhCnLinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2978d00b5a530dcd104d2590d7e7202a99367c54/e2e/$oldId.md"
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $zhCnLinkUrl, "", "", "$newId.md")

$wsZhCn.Columns.Item(9).ColumnWidth = 17.8333333333333
$wsZhCn.Columns.Item(10).ColumnWidth = 20.8333333333333

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("G2").Value = "$newId.$newHash.de-de.xlf"
$wsDeDe.Range("I2").Value = ""
$wsDeDe.Range("I2").Style = "Normal"
$wsDeDe.Range("J2").Value = ""
$wsDeDe.Range("K2").Value = "0001-01-01 00:00:00"

$deDeLinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2978d00b5a530dcd104d2590d7e7202a99367c54/e2e/$oldId.md"
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $deDeLinkUrl, "", "", "$newId.md")

$wsDeDe.Columns.Item(9).ColumnWidth = 17.8333333333333
$wsDeDe.Columns.Item(10).ColumnWidth = 20.8333333333333
